$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 491, shifting all existing data (rows 491:580)
# down to rows 493:582.
$ws.Rows("491:492").Insert()

# Populate the two newly inserted rows with the latest week's data
# (mirrors the structure of the rows that used to occupy 491/492).
$ws.Range("A491").Value = 11
$ws.Range("B491").Value = "Vega Monumental Concepción"
$ws.Range("C491").Value = "Bíobío"
$ws.Range("D491").Value = 45258
$ws.Range("E491").Value = 8
$ws.Range("F491").Value = 100112017
$ws.Range("G491").Value = "Apio"
$ws.Range("H491").Value = "Americana (o)"
$ws.Range("I491").Value = "Primera"
$ws.Range("J491").Value = 330
$ws.Range("K491").Value = 9000
$ws.Range("L491").Value = 10000
$ws.Range("M491").Value = 9545
$ws.Range("N491").Value = "`$/docena de matas"
$ws.Range("O491").Value = "Región de Coquimbo"
$ws.Range("P491").Value = 1591
$ws.Range("Q491").Value = 6
$ws.Range("R491").Value = "Hortaliza"

$ws.Range("A492").Value = 11
$ws.Range("B492").Value = "Vega Monumental Concepción"
$ws.Range("C492").Value = "Bíobío"
$ws.Range("D492").Value = 45258
$ws.Range("E492").Value = 8
$ws.Range("F492").Value = 100112017
$ws.Range("G492").Value = "Apio"
$ws.Range("H492").Value = "Americana (o)"
$ws.Range("I492").Value = "Segunda"
$ws.Range("J492").Value = 220
$ws.Range("K492").Value = 7000
$ws.Range("L492").Value = 8000
$ws.Range("M492").Value = 7545
$ws.Range("N492").Value = "`$/docena de matas"
$ws.Range("O492").Value = "Región de Coquimbo"
$ws.Range("P492").Value = 1258
$ws.Range("Q492").Value = 6
$ws.Range("R492").Value = "Hortaliza"
